# Union County Tax and Fire Rates workbook update:
# Duplicate the "2021-2022" sheet (the first sheet) and place the copy
# before it, then rename the copy to "2022-2023" - this mirrors the
# annual roll-forward of the tax-rate workbook (a new tab for the new
# fiscal year, seeded from last year's numbers).

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("2021-2022")
$source.Copy($source)

$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "2022-2023"

# The new "2022-2023" tab becomes the active tab, cursor left on N8.
$newSheet.Activate()
$newSheet.Range("N8").Select()

# The original "2021-2022" tab's remembered selection moved to B23.
$origSheet = $wb.Worksheets.Item("2021-2022")
$origSheet.Range("B23").Select()

# The previously-active "2020-2021" tab is no longer the active tab;
# its remembered scroll position/selection moved too.
$prevActive = $wb.Worksheets.Item("2020-2021")
$prevActive.Activate()
$prevActive.Application.ActiveWindow.ScrollRow = 4
$prevActive.Range("B10").Select()

# Leave the newly-added "2022-2023" sheet as the active tab, matching
# the workbook's new default view (no explicit activeTab -> first tab).
$newSheet.Activate()
